# Auto-update draw results: append the 2025-10-02 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

# Leading apostrophes force Excel to store the date-like / numeric-like
# strings as literal text (matching the existing rows, which are all text),
# instead of auto-converting "2025-10-02" to a date serial or "251002" to a
# number.
$ws.Cells.Item($row, 1).Value = "'2025-10-02"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251002"
$ws.Cells.Item($row, 4).Value = "6-9-2-9"
$ws.Cells.Item($row, 5).Value = "2025-10-02T21:36:13.045+04:00"

# Restore the default style on the two cells that needed the quote-prefix
# trick, so they don't pick up a "quote prefix" style the other cells (and
# the rest of the table) don't have.
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
